$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cellValues = @{
    "B2" = 1.113956600348445
    "C2" = 0.07041660960562979
    "D2" = 0.002230355565169262
    "E2" = 0.06147272907012846
    "F2" = 4.984343643786985
    "I2" = 3.280589368861371
    "J2" = 0.1775483532458431
    "K2" = 0.9994063224192757
    "L2" = 0.3262083266678957
    "M2" = 0.3088152948878857
    "B3" = 1.1013956312689
    "C3" = 0.06531331347655112
    "D3" = 0.002203931867057918
    "E3" = 0.0617042280321658
    "F3" = 4.944654813428286
    "I3" = 3.255442672115777
    "J3" = 0.1775156476802202
    "K3" = 0.9807796220941611
    "L3" = 0.3260158177424373
    "M3" = 0.3069278121838224
    "B4" = 1.094444579545666
    "C4" = 0.06222300632268229
    "D4" = 0.00219087695088227
    "E4" = 0.0618583785675364
    "F4" = 4.921493858767235
    "I4" = 3.240666928315861
    "J4" = 0.1775099148170582
    "K4" = 0.9700475144179279
    "L4" = 0.3260330249507106
    "M4" = 0.3059494516714132
    "B5" = 1.091803706211067
    "C5" = 0.0609744624101296
    "D5" = 0.002186354025436543
    "E5" = 0.06192422509731399
    "F5" = 4.912359305782971
    "I5" = 3.234812649725157
    "J5" = 0.177511189770005
    "K5" = 0.9658514350865914
    "L5" = 0.3260741512716265
    "M5" = 0.3055962387204367
    "B6" = 1.091376779643269
    "C6" = 0.06076779190254911
    "D6" = 0.002185651141331846
    "E6" = 0.0619353420724349
    "F6" = 4.910860862385761
    "I6" = 3.233850631598102
    "J6" = 0.1775116196713356
    "K6" = 0.9651653941049716
    "L6" = 0.3260830425987322
    "M6" = 0.305540336852939
    "B7" = 1.094408187112862
    "C7" = 0.06220612445582674
    "D7" = 0.002190812725451963
    "E7" = 0.06185925432052652
    "F7" = 4.921369437341411
    "I7" = 3.240587299636843
    "J7" = 0.1775099173859527
    "K7" = 0.9699902064083687
    "L7" = 0.3260334413811279
    "M7" = 0.3059445038861561
    "B8" = 1.109467693689993
    "C8" = 0.06864799846519531
    "D8" = 0.002220586903620259
    "E8" = 0.06155006284905618
    "F8" = 4.970408120037916
    "I8" = 3.27178080248575
    "J8" = 0.1775341001088471
    "K8" = 0.9928376259280469
    "L8" = 0.3261138833635613
    "M8" = 0.3081270573527988
    "B9" = 1.145032888777536
    "C9" = 0.08162646807384988
    "D9" = 0.002304129705258617
    "E9" = 0.06103862800945858
    "F9" = 5.076170379350856
    "I9" = 3.33823325654933
    "J9" = 0.1776952742086131
    "K9" = 1.043232231780621
    "L9" = 0.327343571386649
    "M9" = 0.3138372194715515
    "B10" = 1.174835432528766
    "C10" = 0.09137883455935025
    "D10" = 0.002380874737918148
    "E10" = 0.060720188559511
    "F10" = 5.159752352218788
    "I10" = 3.390296884204375
    "J10" = 0.1778829628794441
    "K10" = 1.083670504081027
    "L10" = 0.3288977550280023
    "M10" = 0.3189020754154193
    "B11" = 1.189190207051666
    "C11" = 0.09586397661169599
    "D11" = 0.002419133039365917
    "E11" = 0.0605876548239106
    "F11" = 5.199059108563858
    "I11" = 3.414690671700527
    "J11" = 0.1779833785011462
    "K11" = 1.102809689351176
    "L11" = 0.3297455607448256
    "M11" = 0.3213946448100273
    "B12" = 1.194740505880844
    "C12" = 0.09756948384924158
    "D12" = 0.002434102115730141
    "E12" = 0.06053923125739669
    "F12" = 5.214128627864511
    "I12" = 3.424030286818606
    "J12" = 0.1780235633239222
    "K12" = 1.110164144120859
    "L12" = 0.3300867968581258
    "M12" = 0.322365577454228
    "B13" = 1.193540062554632
    "C13" = 0.09720185624885858
    "D13" = 0.002430856838955719
    "E13" = 0.06054958181336545
    "F13" = 5.210874911383144
    "I13" = 3.422014285028951
    "J13" = 0.1780148127925472
    "K13" = 1.108575479084521
    "L13" = 0.3300124084665299
    "M13" = 0.3221552676104267
    "B14" = 1.189644540790596
    "C14" = 0.09600414738216045
    "D14" = 0.002420354903135191
    "E14" = 0.06058363568078295
    "F14" = 5.200295180442822
    "I14" = 3.415456998001147
    "J14" = 0.1779866412631108
    "K14" = 1.103412603501454
    "L14" = 0.3297732301177803
    "M14" = 0.3214739822325683
    "B15" = 1.18727332067678
    "C15" = 0.09527144046930403
    "D15" = 0.002413984876394437
    "E15" = 0.06060472414616047
    "F15" = 5.193838873811188
    "I15" = 3.411453786569922
    "J15" = 0.1779696665602089
    "K15" = 1.100264107582774
    "L15" = 0.3296293542118462
    "M15" = 0.321060196272974
    "B16" = 1.173913342438738
    "C16" = 0.09108670694753584
    "D16" = 0.00237844183060254
    "E16" = 0.06072909726295528
    "F16" = 5.157209430303482
    "I16" = 3.388716988366099
    "J16" = 0.1778767028196135
    "K16" = 1.082434670575651
    "L16" = 0.3288451773459187
    "M16" = 0.3187429688033312
    "B17" = 1.165921519225549
    "C17" = 0.08853204357779987
    "D17" = 0.002357494634258828
    "E17" = 0.06080854713328776
    "F17" = 5.135067642617258
    "I17" = 3.374950598697595
    "J17" = 0.177823521697352
    "K17" = 1.071687293696556
    "L17" = 0.3284001320821943
    "M17" = 0.3173696706694251
    "B18" = 1.161399905892836
    "C18" = 0.08706725220160649
    "D18" = 0.002345761342590791
    "E18" = 0.06085540527233047
    "F18" = 5.122453204132142
    "I18" = 3.367099334348893
    "J18" = 0.1777943486420632
    "K18" = 1.065575682862942
    "L18" = 0.3281574098761766
    "M18" = 0.3165975358548607
    "B19" = 1.159881865651442
    "C19" = 0.08657208402692618
    "D19" = 0.002341842742549005
    "E19" = 0.06087147023629136
    "F19" = 5.118202938131418
    "I19" = 3.364452499339222
    "J19" = 0.1777847143031615
    "K19" = 1.063518420578418
    "L19" = 0.3280775067786763
    "M19" = 0.3163391547339671
    "B20" = 1.16676449435343
    "C20" = 0.08880351696581101
    "D20" = 0.002359691898487171
    "E20" = 0.06079996949448674
    "F20" = 5.1374121566908
    "I20" = 3.376409139839751
    "J20" = 0.1778290364560462
    "K20" = 1.072824126505765
    "L20" = 0.328446136368818
    "M20" = 0.3175140239947751
    "B21" = 1.190785645304402
    "C21" = 0.09635575074528901
    "D21" = 0.002423426506801718
    "E21" = 0.06057358542413471
    "F21" = 5.203397684530273
    "I21" = 3.417380256792683
    "J21" = 0.1779948573323615
    "K21" = 1.104926165566752
    "L21" = 0.3298429350771457
    "M21" = 0.3216733585956248
    "B22" = 1.207151905424809
    "C22" = 0.1013328813216958
    "D22" = 0.002467887378369227
    "E22" = 0.06043590893850581
    "F22" = 5.247600868285502
    "I22" = 3.444753121512449
    "J22" = 0.1781158175843984
    "K22" = 1.126529541852051
    "L22" = 0.3308734838416285
    "M22" = 0.3245493689795325
    "B23" = 1.198355958633414
    "C23" = 0.09867268918034711
    "D23" = 0.002443900876363259
    "E23" = 0.06050845172588559
    "F23" = 5.223910129757769
    "I23" = 3.430089135840873
    "J23" = 0.1780501079314831
    "K23" = 1.114942445047717
    "L23" = 0.3303127130674852
    "M23" = 0.3229999834913428
    "B24" = 1.166383157913941
    "C24" = 0.08868077164811439
    "D24" = 0.002358697552026001
    "E24" = 0.06080384376486503
    "F24" = 5.136351843107178
    "I24" = 3.375749536462195
    "J24" = 0.1778265388683451
    "K24" = 1.072309955066657
    "L24" = 0.3284252968960715
    "M24" = 0.3174487076749735
    "B25" = 1.134766212050465
    "C25" = 0.07807777417028205
    "D25" = 0.002278833049848217
    "E25" = 0.06116688177914398
    "F25" = 5.046528571265895
    "I25" = 3.319688466478866
    "J25" = 0.1776395019331325
    "K25" = 1.029000119570213
    "L25" = 0.3268964371138026
    "M25" = 0.3121395989906439
}

foreach ($key in $cellValues.Keys) {
    $ws.Range($key).Value = $cellValues[$key]
}
